$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptos list values (prices, volumes, and the row 31/32 coin swap)
# between ImmutableX and Stellar). Cells whose new value looks like a plain number
# (e.g. "314.43") are first forced to Text format so Excel keeps them as the exact
# original text instead of re-parsing them into a floating point number (which would
# for example turn "7.290" into 7.29 or "0.06570" into 0.0657).
$ws.Cells.Item(2, 4).Value = '29.313.14'
$ws.Cells.Item(2, 5).Value = '  +2.85%  '
$ws.Cells.Item(3, 4).Value = '1.894.78'
$ws.Cells.Item(3, 5).Value = '  +0.89%  '
$ws.Cells.Item(4, 5).Value = '  -0.33%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '314.43'
$ws.Cells.Item(5, 5).Value = '  -0.08%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '1.003'
$ws.Cells.Item(6, 5).Value = '  -0.37%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.5147'
$ws.Cells.Item(7, 5).Value = '  +0.90%  '
$ws.Cells.Item(8, 5).Value = '  -0.38%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.08423'
$ws.Cells.Item(9, 5).Value = '  +0.22%  '
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '42.38'
$ws.Cells.Item(10, 5).Value = '  +1.94%  '
$ws.Cells.Item(11, 5).Value = '  +0.42%  '
$ws.Cells.Item(12, 5).Value = '  +0.32%  '
$ws.Cells.Item(13, 4).Value = '1.896.69'
$ws.Cells.Item(13, 5).Value = '  +1.13%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '20.66'
$ws.Cells.Item(14, 5).Value = '  +0.91%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '7.290'
$ws.Cells.Item(15, 5).Value = '  +0.47%  '
$ws.Cells.Item(16, 5).Value = '  -0.34%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '93.14'
$ws.Cells.Item(17, 5).Value = '  +2.48%  '
$ws.Cells.Item(18, 5).Value = '  -0.03%  '
$ws.Cells.Item(19, 5).Value = '  +0.19%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '17.84'
$ws.Cells.Item(20, 5).Value = '  +0.84%  '
$ws.Cells.Item(21, 5).Value = '  -0.32%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.009'
$ws.Cells.Item(22, 5).Value = '  +1.05%  '
$ws.Cells.Item(23, 4).Value = '29.342.64'
$ws.Cells.Item(23, 5).Value = '  +3.01%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '11.14'
$ws.Cells.Item(24, 5).Value = '  +0.21%  '
$ws.Cells.Item(25, 5).Value = '  -1.91%  '
$ws.Cells.Item(26, 4).Value = '2.109.23'
$ws.Cells.Item(26, 5).Value = '  +0.88%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '159.28'
$ws.Cells.Item(27, 5).Value = '  -1.15%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '20.86'
$ws.Cells.Item(28, 5).Value = '  +0.82%  '
$ws.Cells.Item(29, 5).Value = '  +2.18%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '127.03'
$ws.Cells.Item(30, 5).Value = '  +0.59%  '
$ws.Cells.Item(31, 2).Value = 'Stellar'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.1048'
$ws.Cells.Item(31, 5).Value = '  +0.14%  '
$ws.Cells.Item(32, 2).Value = 'ImmutableX'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.059'
$ws.Cells.Item(32, 5).Value = '  +0.89%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '6.143'
$ws.Cells.Item(33, 5).Value = '  +6.28%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '3.661'
$ws.Cells.Item(34, 5).Value = '  +1.76%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.02483'
$ws.Cells.Item(35, 5).Value = '  +1.66%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.06570'
$ws.Cells.Item(36, 5).Value = '  +1.05%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.2195'
$ws.Cells.Item(37, 5).Value = '  +0.48%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '9.014'
$ws.Cells.Item(38, 5).Value = '  +1.08%  '
$ws.Cells.Item(39, 5).Value = '  +2.38%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.228'
$ws.Cells.Item(40, 5).Value = '  +2.91%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.6521'
$ws.Cells.Item(41, 5).Value = '  +1.31%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.232'
$ws.Cells.Item(42, 5).Value = '  -2.43%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '11.25'
$ws.Cells.Item(43, 5).Value = '  +0.86%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.6061'
$ws.Cells.Item(44, 5).Value = '  -0.04%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '13.19'
$ws.Cells.Item(45, 5).Value = '  +0.94%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '3.673'
$ws.Cells.Item(46, 5).Value = '  -0.54%  '
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '2.048'
$ws.Cells.Item(47, 5).Value = '  +2.04%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.230'
$ws.Cells.Item(48, 5).Value = '  +1.90%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '123.42'
$ws.Cells.Item(49, 5).Value = '  +0.84%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.152'
$ws.Cells.Item(50, 5).Value = '  -3.41%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '77.68'
$ws.Cells.Item(51, 5).Value = '  +0.79%  '
